# Update "想去人数" (interested-count) figures in column F on both the
# "展览" and "全部类型" worksheets, which carry duplicated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> [old value, new value] for column F on each of the two sheets.
$updates = @{
    2  = 150
    4  = 12132
    5  = 1254
    13 = 63
    16 = 353
    17 = 2376
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
